$d = $word.ActiveDocument

# Step 1: Split the "Version" run into "Versi" + "on" while keeping the
# surrounding spellStart/spellEnd proofErr markers in their original places.
$r1 = $d.Range(5, 7)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>on</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# Step 2: Change " 1." to " 2"
$r2 = $d.Range(8, 10)
$r2.Text = "2"

# Step 3: Insert a new run containing "." right after the bookmark
$r3 = $d.Range(9, 9)
$r3.InsertAfter(".")
